$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("展览")
$ws.Range("F3").Value = 928
$ws.Range("F4").Value = 562
$ws.Range("F5").Value = 2323
$ws.Range("F7").Value = 137
$ws.Range("F9").Value = 1202
$ws.Range("F11").Value = 3114
$ws.Range("F14").Value = 1119
$ws.Range("F15").Value = 632
$ws.Range("F16").Value = 13
$ws.Range("F17").Value = 538
$ws.Range("F18").Value = 255
$ws.Range("F19").Value = 11
$ws.Range("F21").Value = 1212
$ws.Range("F22").Value = 1212
$ws.Range("F23").Value = 192
$ws.Range("F26").Value = 214
$ws.Range("F28").Value = 340
$ws.Range("F29").Value = 663
$ws.Range("F31").Value = 159
$ws.Range("F33").Value = 856
$ws.Range("F34").Value = 119
$ws.Range("F36").Value = 319
$ws.Range("F37").Value = 1078
$ws.Range("F38").Value = 5119
$ws.Range("F39").Value = 566
$ws.Range("F40").Value = 305
$ws.Range("F41").Value = 171
$ws.Range("F42").Value = 5
$ws.Range("F44").Value = 16

$ws = $wb.Worksheets("演出")
$ws.Range("F11").Value = 296
$ws.Range("F19").Value = 49
$ws.Range("F23").Value = 408
$ws.Range("F26").Value = 743
$ws.Range("F29").Value = 13

$ws = $wb.Worksheets("本地生活")
$ws.Range("F4").Value = 660

$ws = $wb.Worksheets("全部类型")
$ws.Range("F2").Value = 660
$ws.Range("F6").Value = 928
$ws.Range("F9").Value = 2323
$ws.Range("F11").Value = 137
$ws.Range("F15").Value = 3114
$ws.Range("F17").Value = 296
$ws.Range("F19").Value = 632
$ws.Range("F21").Value = 538
$ws.Range("F22").Value = 255
$ws.Range("F24").Value = 11
$ws.Range("F25").Value = 1212
$ws.Range("F26").Value = 1212
$ws.Range("F30").Value = 214
$ws.Range("F31").Value = 340
$ws.Range("F33").Value = 663
$ws.Range("F35").Value = 408
$ws.Range("F36").Value = 856
$ws.Range("F39").Value = 319
$ws.Range("F40").Value = 1078
$ws.Range("F43").Value = 305
$ws.Range("F44").Value = 171
$ws.Range("F50").Value = 16
